# Weekly update: Agrícola del Norte S.A. de Arica - Betarraga
# Shift existing data rows 167-220 down by two rows (making room for one new
# week of data at the top of the series) and populate the two freed rows
# (167-168) with the new week's values. This mirrors the way each new week
# of observations is inserted just under the header block in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 167 - this pushes old rows 167:220 down to 169:222
# and expands the sheet dimension accordingly.
$ws.Rows("167:168").Insert()

# Row 167 - "Primera" quality, new week
$ws.Cells.Item(167, 1).Value  = 1
$ws.Cells.Item(167, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(167, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(167, 4).Value  = 44559
$ws.Cells.Item(167, 5).Value  = 15
$ws.Cells.Item(167, 6).Value  = 100114014
$ws.Cells.Item(167, 7).Value  = 'Betarraga'
$ws.Cells.Item(167, 8).Value  = 'Sin especificar'
$ws.Cells.Item(167, 9).Value  = 'Primera'
$ws.Cells.Item(167, 10).Value = 800
$ws.Cells.Item(167, 11).Value = 300
$ws.Cells.Item(167, 12).Value = 350
$ws.Cells.Item(167, 13).Value = 325
$ws.Cells.Item(167, 14).Value = '$/paquete 4 unidades'
$ws.Cells.Item(167, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(167, 16).Value = 81
$ws.Cells.Item(167, 17).Value = 4
$ws.Cells.Item(167, 18).Value = 'Hortaliza'

# Row 168 - "Segunda" quality, new week
$ws.Cells.Item(168, 1).Value  = 1
$ws.Cells.Item(168, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(168, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(168, 4).Value  = 44559
$ws.Cells.Item(168, 5).Value  = 15
$ws.Cells.Item(168, 6).Value  = 100114014
$ws.Cells.Item(168, 7).Value  = 'Betarraga'
$ws.Cells.Item(168, 8).Value  = 'Sin especificar'
$ws.Cells.Item(168, 9).Value  = 'Segunda'
$ws.Cells.Item(168, 10).Value = 900
$ws.Cells.Item(168, 11).Value = 300
$ws.Cells.Item(168, 12).Value = 350
$ws.Cells.Item(168, 13).Value = 325
$ws.Cells.Item(168, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(168, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(168, 16).Value = 65
$ws.Cells.Item(168, 17).Value = 5
$ws.Cells.Item(168, 18).Value = 'Hortaliza'
